# Fix typo in column names: "TasP" -> "Task" for the task header row (B1:AB1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AB (2..28) hold the task headers "TasP1".."TasP27" -> "Task1".."Task27"
for ($i = 1; $i -le 27; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "Task$i"
}

# Reflect the cursor move recorded alongside the edit.
$ws.Range("S4").Select()
